# Rename the "_old" / "_new" column header suffixes to the input-file-name
# based suffixes "_FV2410" / "_FV2504" (columns A-J = FV2410 side, L-U = FV2504
# side; column K holds the "diff" header and is left untouched), add an Excel
# Table (ListObject) spanning the used range with those headers, and freeze
# the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns 1-10 (A-J): "<name>_old" -> "<name>_FV2410"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value2 = "$($baseNames[$i])_FV2410"
}

# Columns 12-21 (L-U): "<name>_new" -> "<name>_FV2504"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value2 = "$($baseNames[$i])_FV2504"
}

# Turn the used range into a proper Excel Table ("Table1") with autofilter,
# matching the new header names.
$usedRange = $ws.Range("A1:U64")
$table = $ws.ListObjects.Add(1, $usedRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# Freeze the header row (split below row 1, focus the lower-left pane).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
